# Commit: "added line for icon, didn't work, will circle back later"
# Appends one new log entry to "repair_notes" (row 15) and a matching
# inventory row to "new_inventory" (row 14).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "repair_notes": append a new log row (row 15), modelled on the
# most recent entry (row 14) so the new row picks up the same look/feel
# an end-user would get by duplicating the last line and editing it.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("repair_notes")

$ws1.Range("A14:J14").Copy()
$ws1.Range("A15:J15").PasteSpecial()

$ws1.Cells.Item(15, 1).Value = "7992CL"
$ws1.Cells.Item(15, 2).Value = "Windows-10-10.0.19041-SP0"
$ws1.Cells.Item(15, 3).Value = "Intel(R) Core(TM) i5-8600 CPU @ 3.10GHz"
$ws1.Cells.Item(15, 4).Value = "32 GB"
$ws1.Cells.Item(15, 5).Value = "10.1.3.53"
# F15 (date) keeps the value copied from F14 ("08-24-2021") unchanged.
$ws1.Cells.Item(15, 7).Value = "16:09"
$ws1.Cells.Item(15, 8).Value = "PY_VAR1"
$ws1.Cells.Item(15, 9).Value = "This is me writing in the log"
$ws1.Cells.Item(15, 10).Value = "PY_VAR0"

# ---------------------------------------------------------------------
# Sheet "new_inventory": append a new inventory row (row 14), modelled
# on the previous entry (row 13) for the shared columns.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("new_inventory")

$ws3.Range("A13:Y13").Copy()
$ws3.Range("A14:Y14").PasteSpecial()

$ws3.Cells.Item(14, 1).Value = "7992CL"
# B14 ("Y") keeps the value copied from B13 unchanged.
$ws3.Cells.Item(14, 3).Value = "Dell"
$ws3.Cells.Item(14, 4).Value = "Optiplex"
$ws3.Cells.Item(14, 5).Value = "Tower"
$ws3.Cells.Item(14, 6).Value = "xxxxx"
# H14/I14 ("PY_VAR0"/"PY_VAR1") keep the values copied from row 13 unchanged.
# U14/Y14 ("0") keep the values copied from row 13 unchanged.

# Q14/X14 need text "0"/"1" that differ from what row 13 held, so copy
# the literal text from elsewhere in the column instead of typing a
# bare numeric-looking value (which Excel would coerce to a number).
$ws3.Range("Q11").Copy()
$ws3.Range("Q14").PasteSpecial()
$ws3.Range("X2").Copy()
$ws3.Range("X14").PasteSpecial()
